# Update the "Metadata" worksheet of the ValueSet-age-group workbook.
#
# Summary of the change (FHIR IG ValueSet metadata refresh):
#   - Version bumped from 5.0.0 -> 6.0.0
#   - Date bumped to a newer publish timestamp
#   - Publisher row's (empty) value is now "Alvearie Team"
#   - The duplicated "Contact" / "No display for ContactDetail" row
#     (it appeared twice, rows 10 & 11) is reduced to a single row, whose
#     content is replaced with "Jurisdiction" / "United States of America"
#   - Everything below shifts up by one row, net -1 row overall (A1:B15 -> A1:B14)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Version / Date updates ---
$ws.Range("B3").Value = "6.0.0"
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# --- Publisher row gets its value filled in ---
$ws.Range("B9").Value = "Alvearie Team"

# --- Remove the duplicate "Contact" row (the second copy, row 11) ---
# Everything below (Description, Purpose, Copyright, Immutable) shifts up.
$ws.Rows.Item(11).Delete()

# --- The remaining "Contact" row (now row 10) becomes "Jurisdiction" ---
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

$wb.Save()
